# Auto-generated: update FFXIV market-data value cells per scheduled-runner diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3874.4
$ws.Range("J18").Value = 5998.5
$ws.Range("L18").Value = 5998.5
$ws.Range("N18").Value = -6566.5
$ws.Range("H39").Value = 10658
$ws.Range("I39").Value = 948
$ws.Range("K39").Value = 2844
$ws.Range("M39").Value = -2548
$ws.Range("H70").Value = 2145.4614
$ws.Range("J70").Value = 2448
$ws.Range("L70").Value = 7344
$ws.Range("N70").Value = -7884
$ws.Range("H73").Value = 2145.4614
$ws.Range("J73").Value = 2448
$ws.Range("L73").Value = 7344
$ws.Range("N73").Value = -9216
$ws.Range("H106").Value = 4217.3076
$ws.Range("I106").Value = 2902.0833
$ws.Range("K106").Value = 2902.0833
$ws.Range("M106").Value = -2271.0833
$ws.Range("H113").Value = 6102
$ws.Range("J113").Value = 6286.5
$ws.Range("L113").Value = 6286.5
$ws.Range("N113").Value = -12794.5
$ws.Range("H123").Value = 133497.5
$ws.Range("J123").Value = 133497.5
$ws.Range("L123").Value = 133497.5
$ws.Range("N123").Value = -143297.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 234429.52
$ws.Range("I32").Value = 286486
$ws.Range("K32").Value = 286486
$ws.Range("M32").Value = -286199
$ws.Range("H74").Value = 898876.8
$ws.Range("J74").Value = 2327454
$ws.Range("L74").Value = 2327454
$ws.Range("N74").Value = -2329202
$ws.Range("H77").Value = 898876.8
$ws.Range("J77").Value = 2327454
$ws.Range("L77").Value = 11637270
$ws.Range("N77").Value = -11646006
$ws.Range("H122").Value = 3499.4285
$ws.Range("I122").Value = 3499.4285
$ws.Range("K122").Value = 10498.2855
$ws.Range("M122").Value = -8048.2855

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6873.9165
$ws.Range("I86").Value = 3427.4285
$ws.Range("K86").Value = 3427.4285
$ws.Range("M86").Value = -2304.4285
$ws.Range("H89").Value = 6873.9165
$ws.Range("I89").Value = 3427.4285
$ws.Range("K89").Value = 17137.1425
$ws.Range("M89").Value = -11521.1425

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 303.42856
$ws.Range("I7").Value = 396.4
$ws.Range("J7").Value = 71
$ws.Range("K7").Value = 396.4
$ws.Range("L7").Value = 71
$ws.Range("M7").Value = -283.4
$ws.Range("N7").Value = -297

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 228.35294
$ws.Range("J2").Value = 259.8
$ws.Range("L2").Value = 1558.8
$ws.Range("N2").Value = -1784.8
$ws.Range("H5").Value = 1749.1765
$ws.Range("I5").Value = 1229.625
$ws.Range("K5").Value = 3688.875
$ws.Range("M5").Value = -3576.875
$ws.Range("H25").Value = 3998.8
$ws.Range("I25").Value = 1000
$ws.Range("K25").Value = 3000
$ws.Range("M25").Value = -2831
$ws.Range("H30").Value = 3998.8
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 3000
$ws.Range("M30").Value = -2898
$ws.Range("H129").Value = 8080.2
$ws.Range("J129").Value = 19044.25
$ws.Range("L129").Value = 57132.75
$ws.Range("N129").Value = -67132.75
$ws.Range("H135").Value = 1749.1765
$ws.Range("I135").Value = 1229.625
$ws.Range("K135").Value = 11066.625
$ws.Range("M135").Value = -8531.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1027276.25
$ws.Range("I132").Value = 1218
$ws.Range("J132").Value = 1222716
$ws.Range("K132").Value = 3654
$ws.Range("L132").Value = 3668148
$ws.Range("M132").Value = -1124
$ws.Range("N132").Value = -3673208

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5995.5
$ws.Range("J2").Value = 5995.5
$ws.Range("L2").Value = 5995.5
$ws.Range("N2").Value = -6219.5
$ws.Range("H7").Value = 9485.637000000001
$ws.Range("I7").Value = 4561.3335
$ws.Range("K7").Value = 4561.3335
$ws.Range("M7").Value = -4449.3335
$ws.Range("H22").Value = 4740.8823
$ws.Range("I22").Value = 3333.1667
$ws.Range("K22").Value = 3333.1667
$ws.Range("M22").Value = -3038.1667
$ws.Range("H27").Value = 4740.8823
$ws.Range("I27").Value = 3333.1667
$ws.Range("K27").Value = 3333.1667
$ws.Range("M27").Value = -3226.1667
$ws.Range("H55").Value = 652.2973
$ws.Range("J55").Value = 874.65
$ws.Range("L55").Value = 874.65
$ws.Range("N55").Value = -1220.65
$ws.Range("H122").Value = 3163.5356
$ws.Range("I122").Value = 2891.818
$ws.Range("J122").Value = 4159.8335
$ws.Range("K122").Value = 8675.454000000002
$ws.Range("L122").Value = 12479.5005
$ws.Range("M122").Value = -6225.454000000002
$ws.Range("N122").Value = -17379.5005
$ws.Range("H126").Value = 9485.637000000001
$ws.Range("I126").Value = 4561.3335
$ws.Range("K126").Value = 13684.0005
$ws.Range("M126").Value = -11214.0005
$ws.Range("H130").Value = 37929
$ws.Range("J130").Value = 37929
$ws.Range("L130").Value = 37929
$ws.Range("N130").Value = -47969

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 5150000
$ws.Range("J15").Value = 5150000
$ws.Range("L15").Value = 5150000
$ws.Range("N15").Value = -5150576
$ws.Range("H51").Value = 12535
$ws.Range("I51").Value = 12535
$ws.Range("K51").Value = 12535
$ws.Range("M51").Value = -12025
$ws.Range("H52").Value = 13100
$ws.Range("I52").Value = 13100
$ws.Range("K52").Value = 13100
$ws.Range("M52").Value = -12874
$ws.Range("H81").Value = 104749
$ws.Range("I81").Value = 5928.4287
$ws.Range("J81").Value = 335330.34
$ws.Range("K81").Value = 11856.8574
$ws.Range("L81").Value = 670660.6800000001
$ws.Range("M81").Value = -10795.8574
$ws.Range("N81").Value = -672782.6800000001
$ws.Range("H84").Value = 104749
$ws.Range("I84").Value = 5928.4287
$ws.Range("J84").Value = 335330.34
$ws.Range("K84").Value = 59284.287
$ws.Range("L84").Value = 3353303.4
$ws.Range("M84").Value = -53980.287
$ws.Range("N84").Value = -3363911.4
$ws.Range("H96").Value = 35205
$ws.Range("J96").Value = 100000
$ws.Range("L96").Value = 100000
$ws.Range("N96").Value = -102746
$ws.Range("H122").Value = 3419.8572
$ws.Range("I122").Value = 4155.8
$ws.Range("K122").Value = 12467.4
$ws.Range("M122").Value = -10017.4
$ws.Range("H125").Value = 76357.5
$ws.Range("J125").Value = 76357.5
$ws.Range("L125").Value = 76357.5
$ws.Range("N125").Value = -86197.5
$ws.Range("H136").Value = 60176.59
$ws.Range("I136").Value = 84213.414
$ws.Range("K136").Value = 252640.242
$ws.Range("M136").Value = -250090.242
